$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$sLast = $wb.Worksheets.Add($null, $last)
$sLast.Name = "Sheet7"

$first = $wb.Worksheets.Item(1)
$sFirst = $wb.Worksheets.Add($first)
$sFirst.Name = "sheet8"

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
